# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C271")
$range.Value = 45202
